$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the row for JOSE IVAN GUZMAN ROMERO (old row 21) - shifts everything below up by one
$ws.Rows("21").Delete()

# Re-sort the remaining worker table (now B16:G20) ascending by "Periodo Mora" (column E)
$ws.Range("B16:G20").Sort($ws.Range("E16:E20"))

# The last data row of the table needs the heavier "closing" bottom border that
# used to belong to the (now removed) last row
$ws.Range("B20:J20").Borders.Item(9).LineStyle = 1
$ws.Range("B20:J20").Borders.Item(9).ColorIndex = 1

# Update the summary figures to reflect the updated data set
$ws.Range("E11").Value = 176063
$ws.Range("C13").Value = 4
$ws.Range("F13").Value = 5
